# "Adding banner on main page"
#
# The new banner on the main ("ecofiller") page needs the in-stock/out-of-stock
# flags in column A to be refreshed, and the workbook should open showing that
# sheet (with the view scrolled down a bit on the common_data sheet, and the
# cursor parked lower down on the ecofiller sheet).

$wb = $excel.ActiveWorkbook

# --- Update the stock-flag column (A) on the "ecofiller" sheet -------------
$ecofiller = $wb.Worksheets.Item("ecofiller")

$ecofiller.Range("A1").Value  = 1
$ecofiller.Range("A2").Value  = 2
$ecofiller.Range("A3").Value  = 2
$ecofiller.Range("A4").Value  = 1
$ecofiller.Range("A6").Value  = 1
$ecofiller.Range("A8").Value  = 1
$ecofiller.Range("A9").Value  = 2
$ecofiller.Range("A11").Value = 1
$ecofiller.Range("A14").Value = 2
$ecofiller.Range("A17").Value = 1
$ecofiller.Range("A18").Value = 2
$ecofiller.Range("A22").Value = 2
$ecofiller.Range("A23").Value = 1
$ecofiller.Range("A25").Value = 2
$ecofiller.Range("A29").Value = 2
$ecofiller.Range("A32").Value = 1
$ecofiller.Range("A33").Value = 2
$ecofiller.Range("A38").Value = 1

# --- Scroll the "common_data" sheet down a bit (view only, selection stays) -
$commonData = $wb.Worksheets.Item("common_data")
$commonData.Activate()
$excel.ActiveWindow.ScrollRow = 10

# --- Make "ecofiller" the sheet that is shown/selected when the book opens --
$ecofiller.Activate()
$ecofiller.Range("A12").Select()
